# Fruta / hortaliza, semanal
# Insert a new weekly data point at row 282 (pushing the existing rows
# 282..298 down to 283..299) for the "Zanahoria" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 282, shifting the rest of the
# table (and its formatting) down by one row.
$ws.Rows.Item(282).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(282, 1).Value = 7
$ws.Cells.Item(282, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(282, 3).Value = "Ñuble"
$ws.Cells.Item(282, 4).Value = 44753
$ws.Cells.Item(282, 5).Value = 16
$ws.Cells.Item(282, 6).Value = 100114013
$ws.Cells.Item(282, 7).Value = "Zanahoria"
$ws.Cells.Item(282, 8).Value = "Sin especificar"
$ws.Cells.Item(282, 9).Value = "Primera"
$ws.Cells.Item(282, 10).Value = 120
$ws.Cells.Item(282, 11).Value = 6500
$ws.Cells.Item(282, 12).Value = 7000
$ws.Cells.Item(282, 13).Value = 6750
$ws.Cells.Item(282, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(282, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(282, 16).Value = 338
$ws.Cells.Item(282, 17).Value = 20
$ws.Cells.Item(282, 18).Value = "Hortaliza"
